# Update NATMI TPM-derived values in the LR-pairs sheet (Sema4a-Plxnb2)
# Reflects recomputed ligand/receptor expression and derived specificity/weight
# values after the underlying TPM input data was refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.050905999999999
$ws.Range("H2").Value = 24.152718
$ws.Range("I2").Value = 0.1888708516018927
$ws.Range("J2").Value = 0.1888708516018927
$ws.Range("M2").Value = 11.61289466666667
$ws.Range("N2").Value = 34.838684
$ws.Range("O2").Value = 0.09693042549509606
$ws.Range("P2").Value = 0.09693042549509606
$ws.Range("Q2").Value = 93.49432334923466
$ws.Range("R2").Value = 841.448910143112
$ws.Range("S2").Value = 0.0183073320093926
$ws.Range("T2").Value = 0.0183073320093926
# Row 3
$ws.Range("G3").Value = 8.050905999999999
$ws.Range("H3").Value = 24.152718
$ws.Range("I3").Value = 0.1888708516018927
$ws.Range("J3").Value = 0.1888708516018927
$ws.Range("O3").Value = 0.2981108740043866
$ws.Range("P3").Value = 0.2981108740043866
$ws.Range("Q3").Value = 287.5430939844493
$ws.Range("R3").Value = 2587.887845860044
$ws.Range("S3").Value = 0.05630445464499301
$ws.Range("T3").Value = 0.05630445464499301
# Row 4
$ws.Range("G4").Value = 8.050905999999999
$ws.Range("H4").Value = 24.152718
$ws.Range("I4").Value = 0.1888708516018927
$ws.Range("J4").Value = 0.1888708516018927
$ws.Range("M4").Value = 27.39934733333333
$ws.Range("N4").Value = 82.198042
$ws.Range("O4").Value = 0.2286966748205465
$ws.Range("P4").Value = 0.2286966748205465
$ws.Range("Q4").Value = 220.5895698420173
$ws.Range("R4").Value = 1985.306128578156
$ws.Range("S4").Value = 0.04319413573187775
$ws.Range("T4").Value = 0.04319413573187775
# Row 5
$ws.Range("G5").Value = 8.050905999999999
$ws.Range("H5").Value = 24.152718
$ws.Range("I5").Value = 0.1888708516018927
$ws.Range("J5").Value = 0.1888708516018927
$ws.Range("M5").Value = 45.078635
$ws.Range("N5").Value = 135.235905
$ws.Range("O5").Value = 0.3762620256799708
$ws.Range("P5").Value = 0.3762620256799709
$ws.Range("Q5").Value = 362.9238529933099
$ws.Range("R5").Value = 3266.31467693979
$ws.Range("S5").Value = 0.07106492921562929
$ws.Range("T5").Value = 0.0710649292156293
# Row 6
$ws.Range("I6").Value = 0.2904749299149038
$ws.Range("J6").Value = 0.2904749299149038
$ws.Range("M6").Value = 11.61289466666667
$ws.Range("N6").Value = 34.838684
$ws.Range("O6").Value = 0.09693042549509606
$ws.Range("P6").Value = 0.09693042549509606
$ws.Range("Q6").Value = 143.7900914406538
$ws.Range("R6").Value = 1294.110822965884
$ws.Range("S6").Value = 0.02815585855230984
$ws.Range("T6").Value = 0.02815585855230984
# Row 7
$ws.Range("I7").Value = 0.2904749299149038
$ws.Range("J7").Value = 0.2904749299149038
$ws.Range("O7").Value = 0.2981108740043866
$ws.Range("P7").Value = 0.2981108740043866
$ws.Range("S7").Value = 0.08659373523329492
$ws.Range("T7").Value = 0.08659373523329492
# Row 8
$ws.Range("I8").Value = 0.2904749299149038
$ws.Range("J8").Value = 0.2904749299149038
$ws.Range("M8").Value = 27.39934733333333
$ws.Range("N8").Value = 82.198042
$ws.Range("O8").Value = 0.2286966748205465
$ws.Range("P8").Value = 0.2286966748205465
$ws.Range("Q8").Value = 339.2569011912936
$ws.Range("R8").Value = 3053.312110721642
$ws.Range("S8").Value = 0.0664306505902698
$ws.Range("T8").Value = 0.0664306505902698
# Row 9
$ws.Range("I9").Value = 0.2904749299149038
$ws.Range("J9").Value = 0.2904749299149038
$ws.Range("M9").Value = 45.078635
$ws.Range("N9").Value = 135.235905
$ws.Range("O9").Value = 0.3762620256799708
$ws.Range("P9").Value = 0.3762620256799709
$ws.Range("Q9").Value = 558.1606683538783
$ws.Range("R9").Value = 5023.446015184905
$ws.Range("S9").Value = 0.1092946855390293
$ws.Range("T9").Value = 0.1092946855390293
# Row 10
$ws.Range("G10").Value = 8.230170000000001
$ws.Range("H10").Value = 24.69051
$ws.Range("I10").Value = 0.1930763092661061
$ws.Range("J10").Value = 0.1930763092661061
$ws.Range("M10").Value = 11.61289466666667
$ws.Range("N10").Value = 34.838684
$ws.Range("O10").Value = 0.09693042549509606
$ws.Range("P10").Value = 0.09693042549509606
$ws.Range("Q10").Value = 95.57609729876002
$ws.Range("R10").Value = 860.1848756888402
$ws.Range("S10").Value = 0.01871496881018642
$ws.Range("T10").Value = 0.01871496881018642
# Row 11
$ws.Range("G11").Value = 8.230170000000001
$ws.Range("H11").Value = 24.69051
$ws.Range("I11").Value = 0.1930763092661061
$ws.Range("J11").Value = 0.1930763092661061
$ws.Range("O11").Value = 0.2981108740043866
$ws.Range("P11").Value = 0.2981108740043866
$ws.Range("Q11").Value = 293.9456187686201
$ws.Range("R11").Value = 2645.510568917581
$ws.Range("S11").Value = 0.05755814730486013
$ws.Range("T11").Value = 0.05755814730486013
# Row 12
$ws.Range("G12").Value = 8.230170000000001
$ws.Range("H12").Value = 24.69051
$ws.Range("I12").Value = 0.1930763092661061
$ws.Range("J12").Value = 0.1930763092661061
$ws.Range("M12").Value = 27.39934733333333
$ws.Range("N12").Value = 82.198042
$ws.Range("O12").Value = 0.2286966748205465
$ws.Range("P12").Value = 0.2286966748205465
$ws.Range("Q12").Value = 225.50128644238
$ws.Range("R12").Value = 2029.51157798142
$ws.Range("S12").Value = 0.04415590991578194
$ws.Range("T12").Value = 0.04415590991578194
# Row 13
$ws.Range("G13").Value = 8.230170000000001
$ws.Range("H13").Value = 24.69051
$ws.Range("I13").Value = 0.1930763092661061
$ws.Range("J13").Value = 0.1930763092661061
$ws.Range("M13").Value = 45.078635
$ws.Range("N13").Value = 135.235905
$ws.Range("O13").Value = 0.3762620256799708
$ws.Range("P13").Value = 0.3762620256799709
$ws.Range("Q13").Value = 371.00482941795
$ws.Range("R13").Value = 3339.043464761551
$ws.Range("S13").Value = 0.0726472832352776
$ws.Range("T13").Value = 0.0726472832352776
# Row 14
$ws.Range("G14").Value = 13.96350433333333
$ws.Range("H14").Value = 41.890513
$ws.Range("I14").Value = 0.3275779092170975
$ws.Range("J14").Value = 0.3275779092170975
$ws.Range("M14").Value = 11.61289466666667
$ws.Range("N14").Value = 34.838684
$ws.Range("O14").Value = 0.09693042549509606
$ws.Range("P14").Value = 0.09693042549509606
$ws.Range("Q14").Value = 162.1567050005436
$ws.Range("R14").Value = 1459.410345004892
$ws.Range("S14").Value = 0.03175226612320721
$ws.Range("T14").Value = 0.03175226612320721
# Row 15
$ws.Range("G15").Value = 13.96350433333333
$ws.Range("H15").Value = 41.890513
$ws.Range("I15").Value = 0.3275779092170975
$ws.Range("J15").Value = 0.3275779092170975
$ws.Range("O15").Value = 0.2981108740043866
$ws.Range("P15").Value = 0.2981108740043866
$ws.Range("Q15").Value = 498.7152053286838
$ws.Range("R15").Value = 4488.436847958154
$ws.Range("S15").Value = 0.09765453682123852
$ws.Range("T15").Value = 0.09765453682123852
# Row 16
$ws.Range("G16").Value = 13.96350433333333
$ws.Range("H16").Value = 41.890513
$ws.Range("I16").Value = 0.3275779092170975
$ws.Range("J16").Value = 0.3275779092170975
$ws.Range("M16").Value = 27.39934733333333
$ws.Range("N16").Value = 82.198042
$ws.Range("O16").Value = 0.2286966748205465
$ws.Range("P16").Value = 0.2286966748205465
$ws.Range("Q16").Value = 382.5909052195051
$ws.Range("R16").Value = 3443.318146975546
$ws.Range("S16").Value = 0.07491597858261705
$ws.Range("T16").Value = 0.07491597858261705
# Row 17
$ws.Range("G17").Value = 13.96350433333333
$ws.Range("H17").Value = 41.890513
$ws.Range("I17").Value = 0.3275779092170975
$ws.Range("J17").Value = 0.3275779092170975
$ws.Range("M17").Value = 45.078635
$ws.Range("N17").Value = 135.235905
$ws.Range("O17").Value = 0.3762620256799708
$ws.Range("P17").Value = 0.3762620256799709
$ws.Range("Q17").Value = 629.4557151632516
$ws.Range("R17").Value = 5665.101436469265
$ws.Range("S17").Value = 0.1232551276900347
$ws.Range("T17").Value = 0.1232551276900347
